$d = $word.ActiveDocument
$wns = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

# ---------------------------------------------------------------------------
# 1. Locate the ISP paragraph ("... what it needs.") and the empty spacer
#    paragraph right after it, then:
#      - delete the stray _GoBack bookmark that currently sits at the end
#        of the ISP paragraph (it will be re-created further down, in the
#        middle of the new DIP paragraph)
#      - replace the empty spacer paragraph with one that swaps its
#        jc="left" override for a left indent, and append the whole new
#        "D stands for dependency inversion principle..." paragraph after it
# ---------------------------------------------------------------------------
$i = 0
$ispIdx = -1
foreach ($p in $d.Paragraphs) {
    $i += 1
    if ($p.Range.Text -like "*what it needs*") {
        $ispIdx = $i
    }
}

$ispPara = $d.Paragraphs.Item($ispIdx)
$spacerPara = $d.Paragraphs.Item($ispIdx + 1)

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$startPos = $ispPara.Range.End - 1
$endPos = $spacerPara.Range.End
$mainRange = $d.Range($startPos, $endPos)

$newParasXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="CNormal"/><w:spacing w:after="0"/><w:ind w:left="360"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="CNormal"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">D stands for </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>dependency inversion principle</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t xml:space="preserve">It means that high-level </w:t></w:r><w:r><w:t>modules</w:t></w:r><w:r><w:t xml:space="preserve"> shouldn’t depend on low-level </w:t></w:r><w:r><w:t>modules</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>, instead both should depend on abstractions</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> Also</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> abstractions should not depend on details, details should depend on abstractions.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>
'@

$mainRange.InsertXML($newParasXml)

# ---------------------------------------------------------------------------
# 2. Page-break bookkeeping: the "Bridge" heading and the bullet "&gt;"
#    placeholder paragraph right before the "Visitor" heading both gain a
#    <w:lastRenderedPageBreak/> marker, while the "Visitor" heading loses
#    the one it used to carry.
# ---------------------------------------------------------------------------
$i = 0
$bridgeIdx = -1
$visitorIdx = -1
foreach ($p in $d.Paragraphs) {
    $i += 1
    $t = $p.Range.Text.Trim()
    if ($t -eq "Bridge") {
        $bridgeIdx = $i
    }
    if ($t -eq "Visitor" -and $visitorIdx -eq -1) {
        $visitorIdx = $i
    }
}

$pBridge = $d.Paragraphs.Item($bridgeIdx)
$rBridge = $d.Range($pBridge.Range.Start, $pBridge.Range.End - 1)
$bridgeXml = '<w:p xmlns:w="' + $wns + '"><w:pPr><w:pStyle w:val="CHeading"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Bridge</w:t></w:r></w:p>'
$rBridge.InsertXML($bridgeXml)

$pGt = $d.Paragraphs.Item($visitorIdx - 2)
$rGt = $d.Range($pGt.Range.Start, $pGt.Range.End - 1)
$gtXml = '<w:p xmlns:w="' + $wns + '"><w:pPr><w:pStyle w:val="CNormal"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:spacing w:after="0"/><w:jc w:val="left"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>&gt;</w:t></w:r></w:p>'
$rGt.InsertXML($gtXml)

$pVisitor = $d.Paragraphs.Item($visitorIdx)
$rVisitor = $d.Range($pVisitor.Range.Start, $pVisitor.Range.End - 1)
$visitorXml = '<w:p xmlns:w="' + $wns + '"><w:pPr><w:pStyle w:val="CHeading"/></w:pPr><w:r><w:t>Visitor</w:t></w:r></w:p>'
$rVisitor.InsertXML($visitorXml)
